# Update "Top 50 Cryptocurrencies" sheet with fresh market data
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Top 50 Cryptocurrencies")

$cryptoData = @(
    ,@(2, "Bitcoin", "btc", 98645, 1955265453076, 117514916694, 1.34105)
    ,@(3, "Ethereum", "eth", 3384.98, 408743623550, 57753690385, 7.88246)
    ,@(4, "Tether", "usdt", 0.9997819999999999, 130760277928, 104395714518, -0.11895)
    ,@(5, "Solana", "sol", 260.36, 123643753527, 15001244143, 8.25719)
    ,@(6, "BNB", "bnb", 634.42, 92718719024, 2485956984, 3.93616)
    ,@(7, "XRP", "xrp", 1.39, 79352461642, 17770071804, 24.10278)
    ,@(8, "Dogecoin", "doge", 0.395726, 58339636274, 10048293158, 1.76409)
    ,@(9, "USDC", "usdc", 0.997682, 38285219460, 12956047852, -0.22263)
    ,@(10, "Lido Staked Ether", "steth", 3382.2, 33158319598, 148069164, 7.6953)
    ,@(11, "Cardano", "ada", 0.891133, 31920708422, 1598277482, 11.93)
    ,@(12, "TRON", "trx", 0.200196, 17321666076, 1081275395, 1.34399)
    ,@(13, "Avalanche", "avax", 36.35, 14881045882, 1040305471, 6.49239)
    ,@(14, "Shiba Inu", "shib", 0.00002498, 14733733758, 1620306454, 3.31374)
    ,@(15, "Wrapped stETH", "wsteth", 4016.6, 14529010683, 168190824, 7.87481)
    ,@(16, "Wrapped Bitcoin", "wbtc", 98635, 14437397248, 866759917, 1.49103)
    ,@(17, "Toncoin", "ton", 5.56, 14202717857, 629342951, 3.44887)
    ,@(18, "Sui", "sui", 3.63, 10352847787, 2397340035, 1.18447)
    ,@(19, "Bitcoin Cash", "bch", 494.09, 9808440404, 2115012026, -0.90611)
    ,@(20, "WETH", "weth", 3394.83, 9678336190, 2233267349, 8.261150000000001)
    ,@(21, "Chainlink", "link", 15.27, 9599024349, 1237471114, 4.81535)
    ,@(22, "Pepe", "pepe", 0.00002119, 8959751961, 6872754151, 7.53254)
    ,@(23, "Polkadot", "dot", 6.21, 8947246366, 816079415, 8.464359999999999)
    ,@(24, "Stellar", "xlm", 0.284789, 8574297810, 2311654284, 18.48765)
    ,@(25, "LEO Token", "leo", 8.779999999999999, 8128302174, 3417394, 2.61704)
    ,@(26, "NEAR Protocol", "near", 5.8, 7066502768, 1007516188, 3.93636)
    ,@(27, "Litecoin", "ltc", 90.45, 6818543284, 1443342995, 4.38062)
    ,@(28, "Aptos", "apt", 12.11, 6475051493, 883555179, 3.42643)
    ,@(29, "Wrapped eETH", "weeth", 3573.41, 6229530875, 103976595, 8.2064)
    ,@(30, "Uniswap", "uni", 9.4, 5646687307, 858637673, 6.14347)
    ,@(31, "Cronos", "cro", 0.197491, 5359892708, 120817872, 11.69802)
    ,@(32, "USDS", "usds", 0.99435, 5217775415, 16412763, -0.56931)
    ,@(33, "Hedera", "hbar", 0.133678, 5079959544, 895188955, 5.463)
    ,@(34, "Internet Computer", "icp", 9.66, 4587010509, 272743635, 6.07429)
    ,@(35, "Ethereum Classic", "etc", 27.94, 4201045186, 893965885, 5.76924)
    ,@(36, "Bonk", "bonk", 0.000052, 3925070065, 1720283290, 1.33985)
    ,@(37, "Render", "render", 7.39, 3830491806, 438574850, -0.32458)
    ,@(38, "Kaspa", "kas", 0.150778, 3811511501, 153122650, -0.5954)
    ,@(39, "POL (ex-MATIC)", "pol", 0.470479, 3766297859, 484661119, 6.56518)
    ,@(40, "Bittensor", "tao", 507.51, 3754446505, 285429218, 2.87863)
    ,@(41, "Ethena USDe", "usde", 0.999717, 3685665968, 227695671, -0.3284)
    ,@(42, "WhiteBIT Coin", "wbt", 24.74, 3571018763, 38464423, 2.19327)
    ,@(43, "Dai", "dai", 0.997967, 3436077216, 157037670, -0.25013)
    ,@(44, "dogwifhat", "wif", 3.39, 3398407556, 1279303703, 5.02821)
    ,@(45, "MANTRA", "om", 3.75, 3387528462, 296181418, 2.84473)
    ,@(46, "Artificial Superintelligence Alliance", "fet", 1.28, 3346685530, 490864213, 2.65434)
    ,@(47, "Arbitrum", "arb", 0.7867, 3225646756, 1670400289, 11.73951)
    ,@(48, "Monero", "xmr", 160.85, 2968142680, 86572953, -1.00593)
    ,@(49, "Stacks", "stx", 1.95, 2934163039, 380160497, 0.97964)
    ,@(50, "Filecoin", "fil", 4.69, 2821849071, 588249132, 7.04811)
    ,@(51, "OKB", "okb", 46.63, 2803633244, 19696267, 5.08722)
)

foreach ($row in $cryptoData) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
    $ws1.Cells.Item($r, 6).Value = $row[6]
}

# Update "Top 5 by Market Cap" sheet (mirrors rows 2-6 of sheet1, column D)
$ws2 = $wb.Worksheets.Item("Top 5 by Market Cap")
$ws2.Cells.Item(2, 2).Value = 1955265453076
$ws2.Cells.Item(3, 2).Value = 408743623550
$ws2.Cells.Item(4, 2).Value = 130760277928
$ws2.Cells.Item(5, 2).Value = 123643753527
$ws2.Cells.Item(6, 2).Value = 92718719024

# Update "Summary" sheet
# (leading apostrophe forces the "$..." value to stay text instead of
# being auto-converted to a currency number, matching Excel's own
# quote-prefix behavior; the apostrophe itself is not stored)
$ws3 = $wb.Worksheets.Item("Summary")
$ws3.Cells.Item(2, 2).Value = "'`$4348.48"
$ws3.Cells.Item(3, 2).Value = "XRP (24.10%)"
$ws3.Cells.Item(4, 2).Value = "Monero (-1.01%)"
